# Duplicate Leather data and test file
# Rolls the MemberManagement test data forward: each sheet's test rows are
# replaced with the "next batch" of test users (same shape, incremented
# suffix numbers), mirroring a fresh run of the automated test suite.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: createUser  (A1:F7)  -> UserTest-37..42 / user.testAuto037..042
#                                 become UserTest-43..48 / user.testAuto043..048
#                                 phone numbers 9800000015..20 -> ..21..26
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("createUser")

$ws.Range("A2").Value = "UserTest-43"
$ws.Range("C2").Value = "user.testAuto043"
$ws.Range("D2").Value = 9800000021

$ws.Range("A3").Value = "UserTest-44"
$ws.Range("C3").Value = "user.testAuto044"
$ws.Range("D3").Value = 9800000022

$ws.Range("A4").Value = "UserTest-45"
$ws.Range("C4").Value = "user.testAuto045"
$ws.Range("D4").Value = 9800000023

$ws.Range("A5").Value = "UserTest-46"
$ws.Range("C5").Value = "user.testAuto046"
$ws.Range("D5").Value = 9800000024

$ws.Range("A6").Value = "UserTest-47"
$ws.Range("C6").Value = "user.testAuto047"
$ws.Range("D6").Value = 9800000025

$ws.Range("A7").Value = "UserTest-48"
$ws.Range("C7").Value = "user.testAuto048"
$ws.Range("D7").Value = 9800000026

$ws.Activate() | Out-Null
$ws.Range("A3").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet: editUser  (A1:I6)  -> TestUser-31..35 / userEdit.auto31..35
#                               become TestUser-36..40 / userEdit.auto36..40
#                               phone numbers 9800000010..14 -> ..15..19
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("editUser")

$ws.Range("A2").Value = "TestUser-36"
$ws.Range("C2").Value = "userEdit.auto36"
$ws.Range("D2").Value = 9800000015

$ws.Range("A3").Value = "TestUser-37"
$ws.Range("C3").Value = "userEdit.auto37"
$ws.Range("D3").Value = 9800000016

$ws.Range("A4").Value = "TestUser-38"
$ws.Range("C4").Value = "userEdit.auto38"
$ws.Range("D4").Value = 9800000017

$ws.Range("A5").Value = "TestUser-39"
$ws.Range("C5").Value = "userEdit.auto39"
$ws.Range("D5").Value = 9800000018

$ws.Range("A6").Value = "TestUser-40"
$ws.Range("C6").Value = "userEdit.auto40"
$ws.Range("D6").Value = 9800000019

$ws.Activate() | Out-Null
$ws.Range("A2:A6").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet: resetPassword  (A1:G6) -> TestUser-26..30 / userPass.auto26..30
#                                   become TestUser-35..39 / userPass.auto35..39
#                                   phone numbers unchanged
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("resetPassword")

$ws.Range("A2").Value = "TestUser-35"
$ws.Range("C2").Value = "userPass.auto35"

$ws.Range("A3").Value = "TestUser-36"
$ws.Range("C3").Value = "userPass.auto36"

$ws.Range("A4").Value = "TestUser-37"
$ws.Range("C4").Value = "userPass.auto37"

$ws.Range("A5").Value = "TestUser-38"
$ws.Range("C5").Value = "userPass.auto38"

$ws.Range("A6").Value = "TestUser-39"
$ws.Range("C6").Value = "userPass.auto39"

$ws.Activate() | Out-Null
$ws.Range("A6").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet: checkLogin  (A1:F7) -> UserTest-25..30 / user_logintest.auto25..30
#                                 become UserTest-31..36 / user_logintest.auto31..36
#                                 phone numbers unchanged
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("checkLogin")

$ws.Range("A2").Value = "UserTest-31"
$ws.Range("C2").Value = "user_logintest.auto31"

$ws.Range("A3").Value = "UserTest-32"
$ws.Range("C3").Value = "user_logintest.auto32"

$ws.Range("A4").Value = "UserTest-33"
$ws.Range("C4").Value = "user_logintest.auto33"

$ws.Range("A5").Value = "UserTest-34"
$ws.Range("C5").Value = "user_logintest.auto34"

$ws.Range("A6").Value = "UserTest-35"
$ws.Range("C6").Value = "user_logintest.auto35"

$ws.Range("A7").Value = "UserTest-36"
$ws.Range("C7").Value = "user_logintest.auto36"

# checkLogin stays the active/tab-selected sheet, with C2:C7 selected.
$ws.Activate() | Out-Null
$ws.Range("C2:C7").Select() | Out-Null
